$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking Price values are prefixed with an apostrophe so Excel
# keeps storing them as text (matching the sheet's existing text format),
# consistent with entries like '57.926.91' that use '.' as a thousands separator.

$ws.Range("D2").Value = "'57.926.91"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "'2.448.95"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'522.79"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'131.04"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").Value = "'2.452.12"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "'4.96"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'2.868.75"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "'57.863.61"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'21.73"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "'2.457.05"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "'10.29"
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'313.42"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'64.92"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "'0.403"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'7.20"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "'173.24"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "'3.79"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "'36.24"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("D41").Value = "'0.799"
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "'260.70"
$ws.Range("E44").Value = "  -5.67%  "
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'122.14"
$ws.Range("E47").Value = "  -5.71%  "
$ws.Range("D48").Value = "'0.0496"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "'0.0211"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").Value = "'16.96"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "'16.23"
$ws.Range("E51").Value = "  -3.78%  "
